# "Rename via pdf name instead of most recent"
# Populate the first two Claims rows with data derived from the generated
# PDF claim (name, date range, and billed amount) instead of leaving them
# blank, and restore the previously-active sheet/selection afterwards.

$wb = $excel.ActiveWorkbook

$claims = $wb.Worksheets.Item("Claims")

# Row 2 - McGee, Test
$claims.Range("A2").Value = 45430
$claims.Range("B2").Value = "McGee, Test"
$claims.Range("C2").Value = "5/1/24 - 5/18/24"
$claims.Range("D2").Value = 1235

# Row 3 - Anna, Mary
$claims.Range("A3").Value = 45430
$claims.Range("B3").Value = "Anna, Mary"
$claims.Range("C3").Value = "5/1/24 - 5/18/24"
$claims.Range("D3").Value = 380

# Update the Claims sheet's saved selection to cover the newly-filled rows.
# Selecting a range on a sheet activates that sheet, so re-activate the
# sheet that was active before this edit once the selection is set.
$claims.Range("A2:D18").Select()
$wb.Worksheets.Item("Summary").Activate()
